$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 61
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()

# Row 62
$ws.Range("H62").Value = 3141.0833
$ws.Range("I62").Value = 3119.3
$ws.Range("J62").Value = 3250
$ws.Range("K62").Value = 3119.3
$ws.Range("L62").Value = 3250
$ws.Range("M62").Value = -2495.3
$ws.Range("N62").Value = -4498

# Row 63
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()

# Row 65
$ws.Range("H65").Value = 3141.0833
$ws.Range("I65").Value = 3119.3
$ws.Range("J65").Value = 3250
$ws.Range("K65").Value = 15596.5
$ws.Range("L65").Value = 16250
$ws.Range("M65").Value = -12476.5
$ws.Range("N65").Value = -22490

# Row 66
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()

# Row 129
$ws.Range("H129").Value = 18737.535
$ws.Range("I129").Value = 392.375
$ws.Range("J129").Value = 21795.062
$ws.Range("K129").Value = 1177.125
$ws.Range("L129").Value = 65385.186
$ws.Range("M129").Value = 3822.875
$ws.Range("N129").Value = -75385.186

# Row 137
$ws.Range("H137").Value = 1530.3103
$ws.Range("I137").Value = 1289.95
$ws.Range("J137").Value = 2064.4443
$ws.Range("K137").Value = 3869.85
$ws.Range("L137").Value = 6193.3329
$ws.Range("M137").Value = -1319.85
$ws.Range("N137").Value = -11293.3329

# Row 138
$ws.Range("H138").Value = 1754.6342
$ws.Range("I138").Value = 1181.8334
$ws.Range("J138").Value = 2202.913
$ws.Range("K138").Value = 3545.5002
$ws.Range("L138").Value = 6608.739
$ws.Range("M138").Value = 1594.4998
$ws.Range("N138").Value = -16888.739

$ws = $wb.Worksheets.Item("ARM")
# Row 24
$ws.Range("H24").Value = 30000
$ws.Range("J24").Value = 30000
$ws.Range("L24").Value = 30000
$ws.Range("N24").Value = -30748

# Row 32
$ws.Range("H32").Value = 15731.051
$ws.Range("I32").Value = 16928.928
$ws.Range("J32").Value = 6414.222
$ws.Range("K32").Value = 16928.928
$ws.Range("L32").Value = 6414.222
$ws.Range("M32").Value = -16641.928
$ws.Range("N32").Value = -6988.222

# Row 45
$ws.Range("H45").Value = 1277.8462
$ws.Range("I45").Value = 1314
$ws.Range("J45").Value = 1220
$ws.Range("K45").Value = 1314
$ws.Range("L45").Value = 1220
$ws.Range("M45").Value = -937
$ws.Range("N45").Value = -1974

# Row 100
$ws.Range("H100").Value = 30000
$ws.Range("J100").Value = 30000
$ws.Range("L100").Value = 30000
$ws.Range("N100").Value = -32164

# Row 132
$ws.Range("H132").Value = 5364.12
$ws.Range("I132").Value = 5916.905
$ws.Range("J132").Value = 2462
$ws.Range("K132").Value = 17750.715
$ws.Range("L132").Value = 7386
$ws.Range("M132").Value = -15220.715
$ws.Range("N132").Value = -12446

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 1879.7
$ws.Range("I94").Value = 1566.5
$ws.Range("J94").Value = 2349.5
$ws.Range("K94").Value = 1566.5
$ws.Range("L94").Value = 2349.5
$ws.Range("M94").Value = -1115.5
$ws.Range("N94").Value = -3251.5

# Row 99
$ws.Range("H99").Value = 772.75
$ws.Range("I99").Value = 790
$ws.Range("K99").Value = 790
$ws.Range("M99").Value = 708

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 54.625
$ws.Range("I7").Value = 29.333334
$ws.Range("J7").Value = 87.14286
$ws.Range("K7").Value = 29.333334
$ws.Range("L7").Value = 87.14286
$ws.Range("M7").Value = 83.66666599999999
$ws.Range("N7").Value = -313.14286

# Row 16
$ws.Range("H16").Value = 743.75
$ws.Range("I16").Value = 714.2857
$ws.Range("J16").Value = 950
$ws.Range("K16").Value = 714.2857
$ws.Range("L16").Value = 950
$ws.Range("M16").Value = -427.2857
$ws.Range("N16").Value = -1524

# Row 41
$ws.Range("H41").Value = 4302.273
$ws.Range("J41").Value = 4302.273
$ws.Range("L41").Value = 4302.273
$ws.Range("N41").Value = -5158.273

# Row 59
$ws.Range("H59").Value = 32577.777
$ws.Range("J59").Value = 32577.777
$ws.Range("L59").Value = 32577.777
$ws.Range("N59").Value = -34867.777

# Row 60
$ws.Range("H60").Value = 12930.833
$ws.Range("J60").Value = 13114.782
$ws.Range("L60").Value = 13114.782
$ws.Range("N60").Value = -14136.782

# Row 99
$ws.Range("H99").Value = 2070.087
$ws.Range("I99").Value = 1509.3334
$ws.Range("J99").Value = 2681.818
$ws.Range("K99").Value = 1509.3334
$ws.Range("L99").Value = 2681.818
$ws.Range("M99").Value = -11.33339999999998
$ws.Range("N99").Value = -5677.818

# Row 113
$ws.Range("H113").Value = 743.75
$ws.Range("I113").Value = 714.2857
$ws.Range("J113").Value = 950
$ws.Range("K113").Value = 714.2857
$ws.Range("L113").Value = 950
$ws.Range("M113").Value = 1455.7143
$ws.Range("N113").Value = -5290

# Row 122
$ws.Range("H122").Value = 845.63635
$ws.Range("I122").Value = 811.3333
$ws.Range("K122").Value = 2433.9999
$ws.Range("M122").Value = 16.0001000000002

# Row 126
$ws.Range("H126").Value = 2070.087
$ws.Range("I126").Value = 1509.3334
$ws.Range("J126").Value = 2681.818
$ws.Range("K126").Value = 4528.0002
$ws.Range("L126").Value = 8045.454000000001
$ws.Range("M126").Value = -2058.0002
$ws.Range("N126").Value = -12985.454

$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 89.875
$ws.Range("I12").Value = 2
$ws.Range("J12").Value = 110.15385
$ws.Range("K12").Value = 6
$ws.Range("L12").Value = 330.46155
$ws.Range("M12").Value = 167
$ws.Range("N12").Value = -676.46155

# Row 23
$ws.Range("H23").Value = 86.76470999999999
$ws.Range("I23").Value = 33.166668
$ws.Range("J23").Value = 116
$ws.Range("K23").Value = 99.500004
$ws.Range("L23").Value = 348
$ws.Range("M23").Value = 135.499996
$ws.Range("N23").Value = -818

$ws = $wb.Worksheets.Item("GSM")
# Row 34
$ws.Range("H34").Value = 51515
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()

# Row 76
$ws.Range("H76").Value = 51515
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()

# Row 79
$ws.Range("H79").Value = 51515
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()

# Row 132
$ws.Range("H132").Value = 59161.312
$ws.Range("I132").Value = 70667.14
$ws.Range("J132").Value = 3549.8333
$ws.Range("K132").Value = 212001.42
$ws.Range("L132").Value = 10649.4999
$ws.Range("M132").Value = -209471.42
$ws.Range("N132").Value = -15709.4999

$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 1423.4147
$ws.Range("I132").Value = 1053.7931
$ws.Range("J132").Value = 2316.6667
$ws.Range("K132").Value = 3161.379300000001
$ws.Range("L132").Value = 6950.000100000001
$ws.Range("M132").Value = -631.3793000000005
$ws.Range("N132").Value = -12010.0001

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 1293.2307
$ws.Range("I132").Value = 1164
$ws.Range("J132").Value = 1500
$ws.Range("K132").Value = 3492
$ws.Range("L132").Value = 4500
$ws.Range("M132").Value = -962
$ws.Range("N132").Value = -9560

# Row 136
$ws.Range("H136").Value = 3517.9546
$ws.Range("I136").Value = 4224.1177
$ws.Range("J136").Value = 1117
$ws.Range("K136").Value = 12672.3531
$ws.Range("L136").Value = 3351
$ws.Range("M136").Value = -10122.3531
$ws.Range("N136").Value = -8451
